$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 22839
$ws.Range("I21").Value = 22839
$ws.Range("K21").Value = 22839
$ws.Range("M21").Value = -22371

$ws.Range("H23").Value = 22839
$ws.Range("I23").Value = 22839
$ws.Range("K23").Value = 22839
$ws.Range("M23").Value = -22605

$ws.Range("H29").Value = 1037
$ws.Range("I29").Value = 111
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 333
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -52
$ws.Range("N29").Value = -5062

$ws.Range("H38").Value = 1921.2
$ws.Range("I38").Value = 598.5
$ws.Range("J38").Value = 2803
$ws.Range("K38").Value = 1795.5
$ws.Range("L38").Value = 8409
$ws.Range("M38").Value = -1423.5
$ws.Range("N38").Value = -9153

$ws.Range("H39").Value = 93.64286
$ws.Range("I39").Value = 61.692307
$ws.Range("J39").Value = 509
$ws.Range("K39").Value = 185.076921
$ws.Range("L39").Value = 1527
$ws.Range("M39").Value = 110.923079
$ws.Range("N39").Value = -2119

$ws.Range("H40").Value = 1820.6
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 1867.6666
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 1867.6666
$ws.Range("M40").Value = -1575
$ws.Range("N40").Value = -2217.6666

$ws.Range("H42").Value = 265.2857
$ws.Range("J42").Value = 349.75
$ws.Range("L42").Value = 1049.25
$ws.Range("N42").Value = -1509.25

$ws.Range("H107").Value = 593606.0600000001
$ws.Range("I107").Value = 809172.5600000001
$ws.Range("J107").Value = 798.125
$ws.Range("K107").Value = 809172.5600000001
$ws.Range("L107").Value = 798.125
$ws.Range("M107").Value = -807252.5600000001
$ws.Range("N107").Value = -4638.125

$ws.Range("H132").Value = 208905.25
$ws.Range("I132").Value = 258992.19
$ws.Range("J132").Value = 27821.691
$ws.Range("K132").Value = 776976.5700000001
$ws.Range("L132").Value = 83465.073
$ws.Range("M132").Value = -774446.5700000001
$ws.Range("N132").Value = -88525.073

$ws.Range("H137").Value = 1316.5952
$ws.Range("I137").Value = 783.36365
$ws.Range("J137").Value = 1505.8064
$ws.Range("K137").Value = 2350.09095
$ws.Range("L137").Value = 4517.4192
$ws.Range("M137").Value = 199.9090500000002
$ws.Range("N137").Value = -9617.4192

$ws.Range("H138").Value = 6946786.5
$ws.Range("I138").Value = 2110.276
$ws.Range("J138").Value = 11630406
$ws.Range("K138").Value = 6330.828
$ws.Range("L138").Value = 34891218
$ws.Range("M138").Value = -1190.828
$ws.Range("N138").Value = -34901498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6699.65
$ws.Range("I2").Value = 8082.6875
$ws.Range("K2").Value = 8082.6875
$ws.Range("M2").Value = -7969.6875

$ws.Range("H116").Value = 6699.65
$ws.Range("I116").Value = 8082.6875
$ws.Range("K116").Value = 8082.6875
$ws.Range("M116").Value = -5788.6875

$ws.Range("H122").Value = 1275.826
$ws.Range("I122").Value = 1251
$ws.Range("J122").Value = 1365.2
$ws.Range("K122").Value = 3753
$ws.Range("L122").Value = 4095.6
$ws.Range("M122").Value = -1303
$ws.Range("N122").Value = -8995.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6699.65
$ws.Range("I3").Value = 8082.6875
$ws.Range("K3").Value = 8082.6875
$ws.Range("M3").Value = -7968.6875

$ws.Range("H139").Value = 77999.5
$ws.Range("J139").Value = 77999.5
$ws.Range("L139").Value = 77999.5
$ws.Range("N139").Value = -88279.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H58").Value = 2306.9412
$ws.Range("I58").Value = 1081.3334
$ws.Range("J58").Value = 2975.4546
$ws.Range("K58").Value = 1081.3334
$ws.Range("L58").Value = 2975.4546
$ws.Range("M58").Value = -878.3334
$ws.Range("N58").Value = -3381.4546

$ws.Range("H99").Value = 15626600
$ws.Range("I99").Value = 31250850
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 31250850
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -31249352
$ws.Range("N99").Value = -5346

$ws.Range("H105").Value = 1143.3334
$ws.Range("I105").Value = 1161.25
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1161.25
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 585.75
$ws.Range("N105").Value = -4494

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 15626600
$ws.Range("I126").Value = 31250850
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 93752550
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -93750080
$ws.Range("N126").Value = -11990

$ws.Range("H132").Value = 2249.9565
$ws.Range("I132").Value = 1731.8889
$ws.Range("K132").Value = 5195.6667
$ws.Range("M132").Value = -2665.6667

$ws.Range("H134").Value = 3419.5293
$ws.Range("I134").Value = 993.1429000000001
$ws.Range("J134").Value = 5118
$ws.Range("K134").Value = 2979.4287
$ws.Range("L134").Value = 15354
$ws.Range("M134").Value = -444.4287000000004
$ws.Range("N134").Value = -20424

$ws.Range("H135").Value = 38640
$ws.Range("J135").Value = 38640
$ws.Range("L135").Value = 38640
$ws.Range("N135").Value = -48780

$ws.Range("H136").Value = 2306.9412
$ws.Range("I136").Value = 1081.3334
$ws.Range("J136").Value = 2975.4546
$ws.Range("K136").Value = 3244.0002
$ws.Range("L136").Value = 8926.363799999999
$ws.Range("M136").Value = -694.0001999999999
$ws.Range("N136").Value = -14026.3638

$ws.Range("H137").Value = 40666.668
$ws.Range("J137").Value = 40666.668
$ws.Range("L137").Value = 40666.668
$ws.Range("N137").Value = -50866.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 23
$ws.Range("I12").Value = 12.090909
$ws.Range("J12").Value = 33.909092
$ws.Range("K12").Value = 36.272727
$ws.Range("L12").Value = 101.727276
$ws.Range("M12").Value = 136.727273
$ws.Range("N12").Value = -447.727276

$ws.Range("H33").Value = 214.4375
$ws.Range("I33").Value = 110.7
$ws.Range("J33").Value = 387.33334
$ws.Range("K33").Value = 664.2
$ws.Range("L33").Value = 2324.00004
$ws.Range("M33").Value = -381.2
$ws.Range("N33").Value = -2890.00004

$ws.Range("H107").Value = 263931.75
$ws.Range("J107").Value = 527689.8
$ws.Range("L107").Value = 1583069.4
$ws.Range("N107").Value = -1586909.4

$ws.Range("H131").Value = 2422.1765
$ws.Range("I131").Value = 383.84616
$ws.Range("J131").Value = 2790.2083
$ws.Range("K131").Value = 1151.53848
$ws.Range("L131").Value = 8370.624899999999
$ws.Range("M131").Value = 3888.46152
$ws.Range("N131").Value = -18450.6249

$ws.Range("H140").Value = 7132.421
$ws.Range("I140").Value = 10489.8
$ws.Range("J140").Value = 3402
$ws.Range("K140").Value = 31469.4
$ws.Range("L140").Value = 10206
$ws.Range("M140").Value = -26289.4
$ws.Range("N140").Value = -20566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 688.4706
$ws.Range("I22").Value = 758.3333
$ws.Range("J22").Value = 520.8
$ws.Range("K22").Value = 758.3333
$ws.Range("L22").Value = 520.8
$ws.Range("M22").Value = -463.3333
$ws.Range("N22").Value = -1110.8

$ws.Range("H27").Value = 688.4706
$ws.Range("I27").Value = 758.3333
$ws.Range("J27").Value = 520.8
$ws.Range("K27").Value = 758.3333
$ws.Range("L27").Value = 520.8
$ws.Range("M27").Value = -651.3333
$ws.Range("N27").Value = -734.8

$ws.Range("H46").Value = 501.9
$ws.Range("I46").Value = 533.3333
$ws.Range("K46").Value = 533.3333
$ws.Range("M46").Value = -345.3333

$ws.Range("H55").Value = 465.14285
$ws.Range("I55").Value = 414
$ws.Range("J55").Value = 533.3333
$ws.Range("K55").Value = 414
$ws.Range("L55").Value = 533.3333
$ws.Range("M55").Value = -241
$ws.Range("N55").Value = -879.3333

$ws.Range("H122").Value = 3526.923
$ws.Range("I122").Value = 2712.5
$ws.Range("J122").Value = 3888.889
$ws.Range("K122").Value = 8137.5
$ws.Range("L122").Value = 11666.667
$ws.Range("M122").Value = -5687.5
$ws.Range("N122").Value = -16566.667

$ws.Range("H132").Value = 3534.75
$ws.Range("I132").Value = 2652.8
$ws.Range("J132").Value = 5424.643
$ws.Range("K132").Value = 7958.400000000001
$ws.Range("L132").Value = 16273.929
$ws.Range("M132").Value = -5428.400000000001
$ws.Range("N132").Value = -21333.929

$ws.Range("H136").Value = 4267.6553
$ws.Range("I136").Value = 1358.8889
$ws.Range("J136").Value = 9027.454
$ws.Range("K136").Value = 4076.6667
$ws.Range("L136").Value = 27082.362
$ws.Range("M136").Value = -1526.6667
$ws.Range("N136").Value = -32182.362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8799171
$ws.Range("I136").Value = 10785382
$ws.Range("J136").Value = 3093.2856
$ws.Range("K136").Value = 32356146
$ws.Range("L136").Value = 9279.856800000001
$ws.Range("M136").Value = -32353596
$ws.Range("N136").Value = -14379.8568
